# Fix bug in M.save_excel that dropped the "description" column (and any
# columns beyond "units") when writing the AOCS sheet.
#
# Net effect on the AOCS sheet: the "name"/"value"/"units" header grows a
# new "description" column (inserted before "units", which slides from C
# to D), and the data block is rewritten with the probe/orbiter limits +
# current-type rows; the old "test 1" sample row is preserved, now at the
# bottom (row 9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AOCS")
$ws.Activate()

# Insert a new column for "description" between "value" (B) and the old
# "units" (C) -- "units" slides out to column D, carrying its formatting.
$ws.Columns("C").Insert()

# Header row.
$ws.Cells.Item(1, 1).Value2 = "name"
$ws.Cells.Item(1, 2).Value2 = "value"
$ws.Cells.Item(1, 3).Value2 = "description"
$ws.Cells.Item(1, 4).Value2 = "units"

# New data rows (name, value, description, units).
$data = @(
    ,@("P min temp",   "",         "probe",   "C")
    ,@("P max temp",   "",         "probe",   "C")
    ,@("O min temp",   "",         "orbiter", "C")
    ,@("O max temp",   "",         "orbiter", "C")
    ,@("P max power",  "",         "probe",   "W")
    ,@("O max power",  "",         "orbiter", "W ")
    ,@("current type", "AC or DC", "probe",   "")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $r = $r + 1
}

# Keep the old sample row ("test 1") as the last row.
$ws.Cells.Item(9, 1).Value2 = "test 1"
$ws.Cells.Item(9, 2).Value2 = 3.142857142857143
$ws.Cells.Item(9, 3).Value2 = ""
$ws.Cells.Item(9, 4).Value2 = ""

$ws.Range("A1").Select()
